$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 43 and 44 swap: TheSandbox <-> TrustWalletToken (name, link, price, volume)
$ws.Range("B43").Value = "TrustWalletToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.8480"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -2.52%  "
$ws.Range("B44").Value = "TheSandbox"
$ws.Range("C44").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.4195"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -0.36%  "

# Remaining price / volume updates
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.378.08"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +0.00%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.871.42"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  -0.76%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  +0.20%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "235.44"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -1.06%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.000"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +0.12%  "
$ws.Range("E7").Value = "  -0.42%  "
$ws.Range("E8").Value = "  +0.65%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06559"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -0.27%  "
$ws.Range("E10").Value = "  +7.27%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07939"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +2.62%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "97.59"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -0.95%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.867.39"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -0.99%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.150"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +0.13%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6737"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +0.67%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "282.99"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -0.72%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "30.383.72"
$ws.Range("D17").ClearFormats()
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "5.546"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +3.94%  "
$ws.Range("E19").Value = "  +0.24%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.69"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +0.44%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.111.89"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -0.87%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.000007297"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -0.21%  "
$ws.Range("E23").Value = "  +0.22%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.211"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +0.20%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.293"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -0.05%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "164.83"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -1.34%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "19.15"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +0.42%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.945"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -2.23%  "
$ws.Range("E29").Value = "  -1.51%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.09707"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -1.43%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.446"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -0.60%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.476"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -1.29%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.113"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -1.92%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.04700"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -0.06%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.118"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +1.80%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7049"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -0.62%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.716"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +0.48%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01858"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -0.79%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.329"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -5.25%  "
$ws.Range("E40").Value = "  +0.87%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "73.50"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +1.01%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.944"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -1.22%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.000"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +0.18%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "103.86"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -0.34%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "7.218"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -0.58%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.223"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -2.36%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "937.63"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -5.62%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "34.19"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +0.03%  "
$ws.Range("E51").Value = "  -2.55%  "
